# Add a new tracked LeetCode question row (row 3) to Sheet1:
#   date | Question | URL(hyperlink)
# Mirrors the formatting already used by row 2 (date number format on
# column A, hyperlink style on column C) and widens columns B & C to fit
# the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 values -----------------------------------------------------
# Column A: date (2026-03-12 -> serial 46093), formatted like A2.
$ws.Cells.Item(3, 1).Value = 46093
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

# Column B: plain question text.
$ws.Cells.Item(3, 2).Value = "Stock cooldown DP"

# Column C: hyperlink text + live hyperlink, styled like C2.
$ws.Cells.Item(3, 3).Value = "https://leetcode.com/problems/best-time-to-buy-and-sell-stock-with-cooldown/"
$ws.Hyperlinks.Add($ws.Cells.Item(3, 3), "https://leetcode.com/problems/best-time-to-buy-and-sell-stock-with-cooldown/")
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)   # xlPasteFormats

# --- Column widths ------------------------------------------------------
# Widen column B (new "Question" text) and column C (now-longer URL).
$ws.Columns("B").ColumnWidth = 16.3
$ws.Columns("C").ColumnWidth = 69.3

$ws.Range("A1").Select() | Out-Null
